$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 14.01.2022 15:45"

# D2: change from text "+0.4" to numeric 0.4
$ws.Range("D2").Value = 0.4

# E2: change from text timestamp to a real date/time serial value, matching
# the same date format used by the other rows in column E (style index 2,
# numFmt "YYYY-MM-DD HH:MM:SS") - copy the format from a sibling cell so we
# reuse the existing style instead of creating a new one.
$ws.Range("E2").Value = 44575.64587962963
$ws.Range("E2").NumberFormat = $ws.Range("E3").NumberFormat
